$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Delete row 39 first (gender_n) then row 37 (nick) so row indices
# stay valid while deleting top-down... actually delete bottom row first
# to avoid shifting the row we still need to delete.
$ws.Rows(39).Delete()
$ws.Rows(37).Delete()
